$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data: cluster name -> active cases, in alphabetical order (matches target row order)
$data = @(
    @("126 Racecourse Road Public Housing Tower Flemington", 10),
    @("139 Highett St Apartment Complex Richmond", 10),
    @("3175 The Bays Aged Care Facility Hastings", 17),
    @("3612 BlueCross Glengowrie Outbreak", 10),
    @("3684 Homestyle Aged Care Langford Grange Cranbourne East Outbreak", 20),
    @("4075 Ferndale Gardens Aged Care Services Bayswater North Outbreak", 14),
    @("Australian Lamb Colac East", 14),
    @("Bread Solutions Braeside Outbreak", 20),
    @("CS Square Caroline Springs", 15),
    @("Carton Finishing Pty. Ltd. Campbellfield", 10),
    @("Cedar Meats Australia Brooklyn Outbreak", 11),
    @("Community Kids Bayswater Early Education Centre Bayswater North", 10),
    @("Embracia Aged Care Reservoir Outbreak", 23),
    @("FedEx Station Melbourne Airport", 11),
    @("Green Leaves Early Learning Centre Highlands Craigieburn", 10),
    @("Guardian Childcare Caulfield Outbreak", 15),
    @("Kool Kidz Childcare Narre Warren", 11),
    @("Lantmannen Unibake Australia Mordialloc", 26),
    @("Nido Early School Ascot Vale", 15),
    @("Nido Early School Glenroy", 23),
    @("Northern Health Northern Hospital Epping Emergency Department Tier 1B", 48),
    @("Northern Health The Northern Hospital Epping", 16),
    @("Oceania Meat Processors Laverton North Outbreak", 13),
    @("Pick It Up Fitness Mulgrave Outbreak", 11),
    @("Social Gathering Warrnambool 28 Sep Outbreak", 17),
    @("St Vincents Hospital Emergency Department Melbourne", 37),
    @("Target Distribution Centre Truganina Outbreak", 21),
    @("The Robin Hood Inn Drouin West Outbreak", 16),
    @("The Royal Children's Hospital Melbourne Emergency Department Parkville Tier 1B", 13),
    @("Visy Recycling Springvale", 31),
    @("Werribee Mercy Hospital Emergency Department", 22),
    @("Western Health Sunshine Hospital Emergency Department", 18)
)

# Clear out the old data range below the header (previously went to row 34)
$ws.Range("A2:B34").ClearContents()

# Write new rows starting at row 2
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}
